# Refactored Parser structure. Fixed problems with reading size of classes and methods.
#
# The "Number of Lines" columns on the classNumberOfLines / methodNumberOfLines
# sheets contained bogus "0" placeholders for several classes/methods whose
# line counts were never computed correctly. This script corrects those
# cells to their real values.
#
# The cells hold text (shared-string) values, not numbers, in the original
# workbook. Plainly assigning a numeric-looking string via .Value would make
# Excel coerce it into a real number (or, if forced to text, stamp a new
# "quote-prefixed" cell style onto it). To avoid both side effects we write
# each value through a throw-away text formula and then flatten the formula
# down to a static value via Copy + PasteSpecial(xlPasteValues) - this keeps
# the cell a plain text cell using the default style, just like the rest of
# the sheet.

function Set-TextValue($ws, $cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$wb = $excel.ActiveWorkbook

# --- classNumberOfLines sheet ---
$classSheet = $wb.Worksheets.Item("classNumberOfLines")

# com.zatribune.spring.ecommerce.stock.db.repository.ProductRepository : 0 -> 1
Set-TextValue $classSheet "B4" "1"
# com.zatribune.spring.ecommerce.stock.db.entities.Product : 0 -> 42
Set-TextValue $classSheet "B6" "42"

# --- methodNumberOfLines sheet ---
$methodSheet = $wb.Worksheets.Item("methodNumberOfLines")

# StockApplicationTests() : 0 -> 1
Set-TextValue $methodSheet "C3" "1"
# DevBootstrap.run(java.lang.String[]) : 0 -> 7
Set-TextValue $methodSheet "C11" "7"
# Product.toString() : 0 -> 3
Set-TextValue $methodSheet "C13" "3"
# Product.builder() : 0 -> 3
Set-TextValue $methodSheet "C14" "3"
# Product.getId() : 0 -> 3
Set-TextValue $methodSheet "C15" "3"
# Product.getName() : 0 -> 3
Set-TextValue $methodSheet "C16" "3"
# Product.getAvailableItems() : 0 -> 3
Set-TextValue $methodSheet "C17" "3"
# Product.getReservedItems() : 0 -> 3
Set-TextValue $methodSheet "C18" "3"
# Product.setId(java.lang.Long) : 0 -> 3
Set-TextValue $methodSheet "C19" "3"
# Product.setName(java.lang.String) : 0 -> 3
Set-TextValue $methodSheet "C20" "3"
# Product.setAvailableItems(int) : 0 -> 3
Set-TextValue $methodSheet "C21" "3"
# Product.setReservedItems(int) : 0 -> 3
Set-TextValue $methodSheet "C22" "3"
# Product(java.lang.Long, java.lang.String, int, int) : 0 -> 6
Set-TextValue $methodSheet "C23" "6"
# Product() : 0 -> 2
Set-TextValue $methodSheet "C24" "2"
# StockApplication() : 0 -> 1
Set-TextValue $methodSheet "C28" "1"
